$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 241
$ws.Range("C2").Value = 91200
$ws.Range("D2").Value = 81.69
$ws.Range("F2").Value = 88200

# Row 3
$ws.Range("B3").Value = 793
$ws.Range("C3").Value = 149400
$ws.Range("D3").Value = 77.44
$ws.Range("F3").Value = 123500

# Row 4
$ws.Range("B4").Value = 22
$ws.Range("C4").Value = 62900
$ws.Range("D4").Value = 62.86
$ws.Range("F4").Value = 54100

# Row 5
$ws.Range("C5").Value = 70200
$ws.Range("F5").Value = 51100

# Row 6
$ws.Range("B6").Value = 49
$ws.Range("C6").Value = 72400
$ws.Range("D6").Value = 94.23
$ws.Range("F6").Value = 74400

# Row 7
$ws.Range("B7").Value = 42
$ws.Range("C7").Value = 76500
$ws.Range("D7").Value = 39.25
$ws.Range("F7").Value = 55800

# Row 8
$ws.Range("B8").Value = 9719
$ws.Range("C8").Value = 108500
$ws.Range("D8").Value = 99.51000000000001
$ws.Range("F8").Value = 114100

# Row 9
$ws.Range("B9").Value = 109
$ws.Range("C9").Value = 64100
$ws.Range("D9").Value = 83.84999999999999
$ws.Range("F9").Value = 144800

# Row 10
$ws.Range("B10").Value = 889
$ws.Range("C10").Value = 105800
$ws.Range("D10").Value = 86.73
$ws.Range("F10").Value = 125000
